# Updated questions and excel sheet
#
# - Reorders/relabels the header row (A1:D1): championId, championLevel,
#   championPoints, highestGrade.
# - Adds a "score" / "score/top3" header pair in G1:H1.
# - Adds a grade -> score lookup table in L3:M18 (F..S+ mapped to 1..16).
# - Adjusts column widths for columns A-D and moves the active selection
#   to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "championId"
$ws.Range("B1").Value = "championLevel"
$ws.Range("C1").Value = "championPoints"
$ws.Range("D1").Value = "highestGrade"
$ws.Range("G1").Value = "score"
$ws.Range("H1").Value = "score/top3"

# --- Grade / score lookup table (L3:M18) -----------------------------------
$grades = @("F", "D-", "D", "D+", "C-", "C", "C+", "B-", "B", "B+", "A-", "A", "A+", "S-", "S", "S+")
for ($i = 0; $i -lt $grades.Length; $i++) {
    $row = 3 + $i
    $ws.Range("L$row").Value = $grades[$i]
    $ws.Range("M$row").Value = $i + 1
}

# --- Column widths ----------------------------------------------------------
# Column C keeps its original width (14.33203125) untouched.
$ws.Columns.Item(1).ColumnWidth = 10.166666666666666   # -> stored width 11
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666   # -> stored width 13.5
$ws.Columns.Item(4).ColumnWidth = 11.330729166666666   # -> stored width ~12.1640625 (old column B width)

# --- Selection ---------------------------------------------------------------
$ws.Range("B2").Select()
